$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Simple cell-text replacements (do these first, while row indices are
#     still the original ones, since none of these change the row count) ---
$t.Rows.Item(1).Cells.Item(1).Range.Text  = "0M"
$t.Rows.Item(2).Cells.Item(1).Range.Text  = "0M"
$t.Rows.Item(3).Cells.Item(1).Range.Text  = "0M"
$t.Rows.Item(4).Cells.Item(1).Range.Text  = "69"
$t.Rows.Item(6).Cells.Item(1).Range.Text  = "0.00011"
$t.Rows.Item(7).Cells.Item(1).Range.Text  = "0.00008"
$t.Rows.Item(8).Cells.Item(1).Range.Text  = "0.00001"
$t.Rows.Item(12).Cells.Item(1).Range.Text = "0.00011"

# Collapse the last three multi-run (tab-separated) summary rows down to a
# single value each.
$t.Rows.Item(44).Cells.Item(1).Range.Text = "100"
$t.Rows.Item(45).Cells.Item(1).Range.Text = "0.01"
$t.Rows.Item(46).Cells.Item(1).Range.Text = "231"

# --- Remove the duplicate "0.00007" row (originally row 9) ---
$t.Rows.Item(9).Delete()

# --- Insert a new row (after the row now holding "0.00011", which used to
#     be row 12 "0.00097" and is now row 11 post-delete) and fill it in ---
$t.Rows.Add($t.Rows.Item(12))
$t.Rows.Item(12).Cells.Item(1).Range.Text = "0.00530"
